# Update Nitish Rana's innings-by-innings batting figures (runs, balls,
# fours, sixes) to the corrected values. The sheet stores these numbers
# as text (column format "@", cells are number-stored-as-text), so we
# force the target cells to a text number format before writing the
# value to avoid Excel auto-converting them back to numeric cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = [string]$value
}

# Row 2
Set-TextValue "C2" "58"
Set-TextValue "D2" "35"
Set-TextValue "E2" "4"

# Row 3
Set-TextValue "C3" "26"
Set-TextValue "D3" "13"
Set-TextValue "E3" "6"

# Row 4
Set-TextValue "C4" "0"
Set-TextValue "D4" "1"
Set-TextValue "E4" "0"

# Row 5
Set-TextValue "C5" "29"
Set-TextValue "D5" "20"
Set-TextValue "E5" "3"
Set-TextValue "F5" "1"

# Row 6
Set-TextValue "C6" "2"
Set-TextValue "D6" "4"

# Row 7
Set-TextValue "C7" "0"
Set-TextValue "D7" "1"

# Row 8
Set-TextValue "C8" "9"
Set-TextValue "D8" "10"
Set-TextValue "E8" "1"
Set-TextValue "F8" "0"

# Row 9
Set-TextValue "C9" "24"
Set-TextValue "D9" "18"
Set-TextValue "E9" "2"
Set-TextValue "F9" "1"

# Row 10
Set-TextValue "C10" "87"
Set-TextValue "D10" "61"
Set-TextValue "E10" "10"
Set-TextValue "F10" "4"

# Row 11
Set-TextValue "C11" "5"
Set-TextValue "D11" "6"
Set-TextValue "E11" "1"

# Row 12
Set-TextValue "C12" "9"
Set-TextValue "D12" "14"
Set-TextValue "E12" "1"
Set-TextValue "F12" "0"

# Row 14
Set-TextValue "C14" "81"
Set-TextValue "D14" "53"
Set-TextValue "E14" "13"
Set-TextValue "F14" "1"
